$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Vika 5 (Week 5): Forritun time increased ---
$ws.Range("H19").Value = 120

# --- Vika 7 (Week 7): previously empty week, now filled in ---
# Rannsóknir (row 36)
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 60

# Hönnun (row 38)
$ws.Range("D38").Value = 60
$ws.Range("E38").Value = 60

# Forritun (row 39)
$ws.Range("C39").Value = 180
$ws.Range("D39").Value = 360
$ws.Range("E39").Value = 240

# Prófanir (row 40)
$ws.Range("E40").Value = 60

# --- Vika 8 (Week 8): previously empty week, now filled in ---
# Forritun (row 49)
$ws.Range("C49").Value = 150
$ws.Range("D49").Value = 240
$ws.Range("E49").Value = 180

# Prófanir (row 50)
$ws.Range("D50").Value = 60

# Frágangur (row 51)
$ws.Range("D51").Value = 60
$ws.Range("E51").Value = 120

# --- Samantekt (Summary) total formula updated to sum weekly totals directly,
#     and a new "hours" conversion cell added next to it ---
$ws.Range("D62").Formula = "=SUM(J12,J22,J32,J42,J52)"
$ws.Range("F62").Formula = "=D62/60"

# --- Update the active selection to reflect where the editor left off ---
$ws.Range("J60").Select()
